$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 3 for columns A, Q, R, Z, AB
$colsToSwap = @("A", "Q", "R", "Z", "AB")

foreach ($col in $colsToSwap) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $val2 = $cell2.Value2
    $val3 = $cell3.Value2
    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
